$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying H1's format (bold/border/centered)
# so the same shared cell style is reused rather than creating a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-10
$iValues = @(6, 6, 7, 8, 9, 3, 8, 6, 7)
$jValues = @(7, 7, 8, 8, 9, 3, 8, 6, 7)

for ($r = 0; $r -lt 9; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
